$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new confirmed "Weekend Getaway" bookings for SNOW-145448, logged 2026-02-17.
# (7 total bookings after this update.)

# ---- Row 7 : 2026-02-21 ----
$ws.Range("A7").Value = "SNOW-145448"
$ws.Range("C7").Value = "a"
$ws.Range("D7").Value = "a@b.com"
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = "Weekend Getaway"
$ws.Range("H7").Value = 12800
$ws.Range("I7").Value = 25600
$ws.Range("J7").Value = "Confirmed"

# Numeric-looking / date-looking values have to be forced to text (NumberFormat "@"),
# otherwise Excel auto-converts them to a number/date serial. Re-apply the plain
# (General, unstyled) format afterwards by pasting it from an already-plain cell so
# the stored value stays textual without leaving a lingering custom cell style.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2026-02-21"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "12"
$ws.Range("A7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "2026-02-17"
$ws.Range("A7").Copy()
$ws.Range("K7").PasteSpecial(-4122)

# Special Requests left blank (stored as an empty string, same as rows 5-6).
$ws.Range("L7").Formula = '=""'

# ---- Row 8 : 2026-02-22 ----
$ws.Range("A8").Value = "SNOW-145448"
$ws.Range("C8").Value = "a"
$ws.Range("D8").Value = "a@b.com"
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = "Weekend Getaway"
$ws.Range("H8").Value = 12800
$ws.Range("I8").Value = 25600
$ws.Range("J8").Value = "Confirmed"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2026-02-22"
$ws.Range("A8").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "12"
$ws.Range("A8").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "2026-02-17"
$ws.Range("A8").Copy()
$ws.Range("K8").PasteSpecial(-4122)

$ws.Range("L8").Formula = '=""'
